# Update the "Förändrad" (Changed) date column C for rows 2-10 from
# serial date 45208 (2023-10-09) to 45212 (2023-10-13), preserving the
# existing cell formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}
